# Insert a new weekly price-report row before row 893 (pushing existing
# rows 893-996 down to 894-997) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(893).Insert()

$ws.Cells.Item(893, 1).Value  = 8
$ws.Cells.Item(893, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(893, 3).Value  = "Coquimbo"
$ws.Cells.Item(893, 4).Value  = 44946
$ws.Cells.Item(893, 5).Value  = 4
$ws.Cells.Item(893, 6).Value  = 100112004
$ws.Cells.Item(893, 7).Value  = "Cebolla"
$ws.Cells.Item(893, 8).Value  = "Sin especificar"
$ws.Cells.Item(893, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(893, 10).Value = 2000
$ws.Cells.Item(893, 11).Value = 10000
$ws.Cells.Item(893, 12).Value = 11000
$ws.Cells.Item(893, 13).Value = 10500
$ws.Cells.Item(893, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(893, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(893, 16).Value = 656
$ws.Cells.Item(893, 17).Value = 16
$ws.Cells.Item(893, 18).Value = "Hortaliza"
